# scenarios.xlsx - update the latest test-run results on the "Plan1" sheet
# (mirrors a later automated test pass: a couple of run timestamps were
# refreshed, one more test now reports "No"/"Failed", and the active
# selection moved down to B10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Row 2 (CT 01): refresh the vOutData run timestamp
$ws.Range("H2").Value = "28_04_2020--21_27_33 556"

# Row 5 (CT 04): vCreateMovement flag flips from Yes to No, and the
# vOutData run timestamp is refreshed
$ws.Range("B5").Value = "No"
$ws.Range("H5").Value = "28_04_2020--21_22_15 376"

# Row 6 (CT 05): Status flips from Passed to Failed, and the vOutData
# run timestamp is refreshed
$ws.Range("C6").Value = "Failed"
$ws.Range("H6").Value = "28_04_2020--21_28_27 839"

# Move/save the active selection to B10, as recorded in the sheet view
$ws.Range("B10").Select()
